$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DPLKINV093-001")
$ws2 = $wb.Worksheets.Item("DPLKINV093-002")

# The Bank ID generated during the approval run turned out to be PAU0228
# (was previously recorded as PAU0226 / "Sesuai hasil generate").

# --- DPLKINV093-002 (sheet2) explanation text first ---
$ws2.Range("F2").Value = "Username : 32070;`nPassword : bni1234;`nRole : 18/19 - Pimpinan Kelompok Investasi/Pengelola Investasi;`nBank ID : PAU0228;`nStatus Verifikasi : 0 : Kembalikan ke Register;`nKeterangan Verifikasi : DATA AKAN DIKEMBALIKAN UNTUK DIEDIT"

# --- Bank ID cells on both sheets ---
$ws1.Range("M2").Value = "PAU0228"
$ws2.Range("M2").Value = "PAU0228"

# --- DPLKINV093-001 (sheet1) explanation text ---
$ws1.Range("F2").Value = "Username : 32070;`nPassword : bni1234;`nRole : 18/19 - Pimpinan Kelompok Investasi/Pengelola Investasi;`nBank ID : PAU0228;`nStatus Verifikasi : 1 : Setuju;`nKeterangan Verifikasi : DATA APPROVAL"

# --- View/selection state ---
# Sheet "DPLKINV093-002" stays the active tab throughout (tabSelected unaffected);
# switch to sheet1 briefly to move its saved selection/scroll, then switch back.
$ws1.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws1.Range("G2").Select()

$ws2.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws2.Range("O2").Select()
